$d = $word.ActiveDocument

$pairs = @(
    @("23×79=", "51×39="),
    @("15×83=", "36×47="),
    @("78×74=", "82×31="),
    @("57×13=", "98×40="),
    @("93×63=", "60×86="),
    @("80×24=", "41×46="),
    @("76×76=", "26×72="),
    @("86×92=", "13×65="),
    @("34×16=", "41×25="),
    @("78×33=", "14×16="),
    @("15×91=", "27×30="),
    @("96×16=", "96×24="),
    @("50×79=", "96×91="),
    @("80×45=", "46×68="),
    @("31×14=", "98×92="),
    @("15×15=", "99×26="),
    @("21×13=", "73×32="),
    @("41×33=", "40×89="),
    @("15×82=", "19×67="),
    @("32×89=", "14×81="),
    @("52×77=", "47×96="),
    @("35×40=", "57×63="),
    @("28×89=", "39×68="),
    @("15×88=", "41×74="),
    @("60×46=", "24×55=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
